# Incremento do tamanho dos telefones no layout
# Update the phone-related field sizes (column C) on the "Cliente" sheet
# from 10 to 11. Column D recalculates automatically because it holds a
# shared SUM formula that chains down the rows.

$wb = $excel.ActiveWorkbook

$wsCliente = $wb.Worksheets.Item("Cliente")
$wsImovel  = $wb.Worksheets.Item("Imovel")

# telefone_usuario / celular_usuario
$wsCliente.Range("C17").Value = 11
$wsCliente.Range("C18").Value = 11

# telefone_proprietario / celular_proprietario
$wsCliente.Range("C27").Value = 11
$wsCliente.Range("C28").Value = 11

# telefone_responsavel / celular_responsavel
$wsCliente.Range("C44").Value = 11
$wsCliente.Range("C45").Value = 11

# Reflect the view/selection state captured in the workbook when it was
# saved: the "Imovel" sheet's selection moved while the "Cliente" sheet
# became the active (selected) tab.
$wsImovel.Range("E41").Select()

$wsCliente.Activate()
$wsCliente.Range("B18").Select()
